$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (old "Student" row) entirely - final table is A1:E3
$ws.Rows("4").Delete()

# Make sure the whole target range is stored as text (matches source
# workbook, which keeps even numeric-looking values as text cells)
$ws.Range("A1:E3").NumberFormat = "@"

# Header row (A1 stays blank, same as source)
$ws.Range("B1").Value = "GK1"
$ws.Range("C1").Value = "CK1"
$ws.Range("D1").Value = "Bonu"
$ws.Range("E1").Value = "Total"

# Row 2
$ws.Range("A2").Value = "Hoa Pham 2"
$ws.Range("B2").Value = "100"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "60"
$ws.Range("E2").Value = "161"

# Row 3
$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = "60"
$ws.Range("C3").Value = "2"
$ws.Range("D3").Value = "300"
$ws.Range("E3").Value = "362"
